# Auto-generated script applying cell-level value updates to match target diff
# (cryptos.xlsx price/volume refresh, commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column D holds price text that often LOOKS numeric (e.g. '1.001', '24.521.67').
# Pre-set those cells to Text format so Excel stores the literal string instead of
# silently coercing it to a floating point number (matches source data being text).
$priceCells = @('D2', 'D3', 'D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D50', 'D51')
foreach ($addr in $priceCells) { $ws.Range($addr).NumberFormat = '@' }

$ws.Range('D2').Value = '24.521.67'
$ws.Range('E2').Value = '  -0.45%  '
$ws.Range('D3').Value = '1.697.26'
$ws.Range('E3').Value = '  +0.07%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '316.85'
$ws.Range('E5').Value = '  +0.50%  '
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').Value = '0.3920'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '0.4054'
$ws.Range('E8').Value = '  +0.38%  '
$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').Value = '54.36'
$ws.Range('E9').Value = '  +2.59%  '
$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D10').Value = '1.484'
$ws.Range('E10').Value = '  -2.23%  '
$ws.Range('B11').Value = 'BinanceUSD'
$ws.Range('C11').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D11').Value = '1.003'
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('D12').Value = '0.08805'
$ws.Range('E12').Value = '  -0.28%  '
$ws.Range('D13').Value = '26.13'
$ws.Range('E13').Value = '  +10.92%  '
$ws.Range('D14').Value = '7.450'
$ws.Range('E14').Value = '  +0.66%  '
$ws.Range('D15').Value = '8.097'
$ws.Range('E15').Value = '  -0.44%  '
$ws.Range('D16').Value = '0.00001357'
$ws.Range('E16').Value = '  +2.82%  '
$ws.Range('D17').Value = '1.697.76'
$ws.Range('E17').Value = '  -0.18%  '
$ws.Range('D18').Value = '97.61'
$ws.Range('E18').Value = '  -1.84%  '
$ws.Range('D19').Value = '0.07185'
$ws.Range('E19').Value = '  +2.43%  '
$ws.Range('D20').Value = '20.45'
$ws.Range('E20').Value = '  +3.77%  '
$ws.Range('D21').Value = '7.298'
$ws.Range('E21').Value = '  +3.40%  '
$ws.Range('D22').Value = '1.002'
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('D23').Value = '14.34'
$ws.Range('E23').Value = '  -1.84%  '
$ws.Range('D24').Value = '24.527.05'
$ws.Range('E24').Value = '  -0.42%  '
$ws.Range('D25').Value = '3.016'
$ws.Range('E25').Value = '  -3.30%  '
$ws.Range('D26').Value = '2.343'
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('D27').Value = '22.95'
$ws.Range('E27').Value = '  +1.49%  '
$ws.Range('D28').Value = '168.89'
$ws.Range('E28').Value = '  +3.54%  '
$ws.Range('D29').Value = '5.907'
$ws.Range('E29').Value = '  +14.83%  '
$ws.Range('D30').Value = '144.69'
$ws.Range('E30').Value = '  +6.73%  '
$ws.Range('D31').Value = '8.398'
$ws.Range('E31').Value = '  -4.09%  '
$ws.Range('D32').Value = '1.881.13'
$ws.Range('E32').Value = '  -0.42%  '
$ws.Range('D33').Value = '2.179'
$ws.Range('E33').Value = '  +11.23%  '
$ws.Range('D34').Value = '0.08775'
$ws.Range('E34').Value = '  -1.89%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = '1.051'
$ws.Range('E35').Value = '  -1.20%  '
$ws.Range('B36').Value = 'InternetComputer(DFINITY)'
$ws.Range('C36').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D36').Value = '7.192'
$ws.Range('E36').Value = '  -5.01%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '0.03100'
$ws.Range('E37').Value = '  +7.56%  '
$ws.Range('D38').Value = '0.2807'
$ws.Range('E38').Value = '  +2.27%  '
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').Value = '0.8586'
$ws.Range('E39').Value = '  +12.46%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '10.88'
$ws.Range('E40').Value = '  -1.25%  '
$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').Value = '0.09180'
$ws.Range('E41').Value = '  +0.35%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').Value = '14.20'
$ws.Range('E42').Value = '  -1.38%  '
$ws.Range('D43').Value = '1.482'
$ws.Range('E43').Value = '  +1.69%  '
$ws.Range('D44').Value = '17.42'
$ws.Range('E44').Value = '  +10.06%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = '0.7499'
$ws.Range('E45').Value = '  +4.65%  '
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').Value = '2.662'
$ws.Range('E46').Value = '  +3.34%  '
$ws.Range('D47').Value = '4.261'
$ws.Range('E47').Value = '  +1.37%  '
$ws.Range('D48').Value = '1.388'
$ws.Range('E48').Value = '  +3.88%  '
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('D50').Value = '140.30'
$ws.Range('E50').Value = '  +0.30%  '
$ws.Range('D51').Value = '0.08226'
$ws.Range('E51').Value = '  +3.27%  '
